# Apply hybrid bold + color (#2C3E50) highlighting to quantitative impact
# metrics (percentages, dollar amounts, large counts) in specific bullet /
# impact paragraphs, matching the "quantitative metrics highlighting"
# commit.

$d = $word.ActiveDocument
$metricColor = 5258796  # wdColor value for RGB 2C3E50 (44, 62, 80)

function Format-MetricInRange {
    param(
        $Document,
        $SearchStart,
        $SearchEnd,
        [string]$MetricText
    )

    $r = $Document.Range($SearchStart, $SearchEnd)
    $find = $r.Find
    $find.ClearFormatting()
    $find.Text = $MetricText
    $find.MatchCase = $true
    $find.MatchWholeWord = $false
    $find.MatchWildcards = $false
    $find.Forward = $true
    $find.Wrap = 0
    $result = $find.Execute()
    if ($result) {
        $r.Font.Bold = $true
        $r.Font.Color = $metricColor
        return $r.End
    }
    return $SearchStart
}

function Format-MetricsInParagraph {
    param(
        $Document,
        [string]$ParagraphText,
        [string[]]$Metrics
    )

    for ($i = 1; $i -le $Document.Paragraphs.Count; $i++) {
        $p = $Document.Paragraphs.Item($i)
        $t = $p.Range.Text.Trim()
        if ($t.EndsWith($ParagraphText)) {
            $paraStart = $p.Range.Start
            $paraEnd = $p.Range.End
            $cursor = $paraStart
            foreach ($metric in $Metrics) {
                $cursor = Format-MetricInRange $Document $cursor $paraEnd $metric
            }
            return $true
        }
    }
    return $false
}

# 1. "• Discovered systematic race coding errors ... from 23% to 64%"
Format-MetricsInParagraph $d "Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%" @("23%", "64%") | Out-Null

# 2. "• Achieved 87% prediction accuracy ... reducing polling error margins from ±4.2% to ±2.1%"
$plusMinus = [char]0x00B1
$para2Text = "Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from " + $plusMinus + "4.2% to " + $plusMinus + "2.1%"
Format-MetricsInParagraph $d $para2Text @("87%", "71%", ($plusMinus + "4.2%"), ($plusMinus + "2.1%")) | Out-Null

# 3. "• Wrote RFP and analyzed bids from 1,200 vendors for research platform development"
Format-MetricsInParagraph $d "Wrote RFP and analyzed bids from 1,200 vendors for research platform development" @("1,200") | Out-Null

# 4. "• Created comprehensive meta-analysis framework ... $400M ... now valued at $1B+"
Format-MetricsInParagraph $d "Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+" @("`$400M", "`$1B") | Out-Null

# 5. "• Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M"
Format-MetricsInParagraph $d "Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M" @("73.5%", "`$4.7M") | Out-Null

# 6. "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%" (short form, no polling margin clause)
Format-MetricsInParagraph $d "Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%" @("87%", "71%") | Out-Null

Write-Output "done"
